$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 35) - Prediction / Error updated
$ws.Range("D2").Value = 0.999999089123995
$ws.Range("E2").Value = 0.999999089123995

# Row 3 (Control 1) - Success flips TRUE -> FALSE, Prediction / Error updated
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.9978621195641789
$ws.Range("E3").Value = 0.9978621195641789

# Row 4 (Control 31) - Success flips TRUE -> FALSE, Prediction / Error updated
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.9999998599473182
$ws.Range("E4").Value = 0.9999998599473182

# Row 5 (Control 14) - Success flips FALSE -> TRUE, Prediction / Error updated
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.009823860515630559
$ws.Range("E5").Value = 0.009823860515630559

# Row 6 (Control 19) - Prediction / Error updated
$ws.Range("D6").Value = 0.00095296686483922
$ws.Range("E6").Value = 0.00095296686483922

# Row 7 (MDD 41) - Prediction / Error updated
$ws.Range("D7").Value = (2.46980613245879 * [math]::Pow(10,-17))
$ws.Range("E7").Value = 1

# Row 8 (MDD 8) - Prediction / Error updated
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0

# Row 9 (MDD 15) - Prediction / Error updated
$ws.Range("D9").Value = 0.001559825461911348
$ws.Range("E9").Value = 0.9984401745380886

# Row 10 (MDD 16) - Prediction / Error updated
$ws.Range("D10").Value = 0.9999999999727156
$ws.Range("E10").Value = (2.7284396963978 * [math]::Pow(10,-11))

# Row 11 (MDD 33) - Prediction / Error / Cross Entropy Loss / Success % updated
$ws.Range("D11").Value = 0.9999999999999991
$ws.Range("E11").Value = (8.881784197001252 * [math]::Pow(10,-16))
$ws.Range("F11").Value = 8.05518627166748
$ws.Range("G11").Value = 0.5
